$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ink / Varnish" rows 2-9 describe two process groups:
#   rows 2-5 = "Flat  2p" @ 0.07, cycling through Yellow/Cyan/Magenta/Black
#   rows 6-9 = "Banner"   @ 0.08, cycling through Cyan/Yellow/Black/Magenta
# The edit swaps which process group owns which row block (and which color
# goes with which row), i.e. a permutation of the (Process, Color, Quantity)
# triple living in columns B, D, E across rows 2-9:
#   new row 2 <- old row 6      new row 6 <- old row 5
#   new row 3 <- old row 9      new row 7 <- old row 4
#   new row 4 <- old row 8      new row 8 <- old row 2
#   new row 5 <- old row 7      new row 9 <- old row 3
#
# Columns B/D/E are shared-string ("text") cells (e.g. "0.07"/"0.08" are
# stored as text, not numbers) - plain .Value assignment would coerce the
# numeric-looking strings into real numbers, changing their cell type.
# Using Copy + PasteSpecial(values-only) instead preserves the original
# text type and leaves destination formatting untouched.

$srcRows = 2..9
foreach ($r in $srcRows) {
    $stageRow = 100 + $r
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("B$stageRow").PasteSpecial(-4163) | Out-Null
    $ws.Range("D$r").Copy() | Out-Null
    $ws.Range("D$stageRow").PasteSpecial(-4163) | Out-Null
    $ws.Range("E$r").Copy() | Out-Null
    $ws.Range("E$stageRow").PasteSpecial(-4163) | Out-Null
}

$rowMap = @{2=106; 3=109; 4=108; 5=107; 6=105; 7=104; 8=102; 9=103}

foreach ($dstRow in $rowMap.Keys) {
    $stageRow = $rowMap[$dstRow]
    $ws.Range("B$stageRow").Copy() | Out-Null
    $ws.Range("B$dstRow").PasteSpecial(-4163) | Out-Null
    $ws.Range("D$stageRow").Copy() | Out-Null
    $ws.Range("D$dstRow").PasteSpecial(-4163) | Out-Null
    $ws.Range("E$stageRow").Copy() | Out-Null
    $ws.Range("E$dstRow").PasteSpecial(-4163) | Out-Null
}

# Remove the scratch staging rows so they don't show up in the sheet.
$ws.Range("A102:J109").Clear() | Out-Null

$wb.Save()
